$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 198; this shifts the existing
# rows 198-248 down to rows 200-250 (matching the rest of the diff, which is
# just that shift).
$ws.Range("A198:A199").EntireRow.Insert()

# Populate the two newly inserted rows (198 and 199) with the new weekly data.
$ws.Range("A198").Value = 7
$ws.Range("B198").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C198").Value = "Ñuble"
$ws.Range("D198").Value = 44641
$ws.Range("E198").Value = 16
$ws.Range("F198").Value = 100112002
$ws.Range("G198").Value = "Pimiento"
$ws.Range("H198").Value = "Zafiro rojo"
$ws.Range("I198").Value = "Primera"
$ws.Range("J198").Value = 120
$ws.Range("K198").Value = 16000
$ws.Range("L198").Value = 17000
$ws.Range("M198").Value = 16500
$ws.Range("N198").Value = "`$/caja 15 kilos"
$ws.Range("O198").Value = "Región del Maule"
$ws.Range("P198").Value = 1100
$ws.Range("Q198").Value = 15
$ws.Range("R198").Value = "Hortaliza"

$ws.Range("A199").Value = 7
$ws.Range("B199").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C199").Value = "Ñuble"
$ws.Range("D199").Value = 44641
$ws.Range("E199").Value = 16
$ws.Range("F199").Value = 100112002
$ws.Range("G199").Value = "Pimiento"
$ws.Range("H199").Value = "Zafiro verde"
$ws.Range("I199").Value = "Primera"
$ws.Range("J199").Value = 120
$ws.Range("K199").Value = 10000
$ws.Range("L199").Value = 11000
$ws.Range("M199").Value = 10500
$ws.Range("N199").Value = "`$/caja 15 kilos"
$ws.Range("O199").Value = "Región del Maule"
$ws.Range("P199").Value = 700
$ws.Range("Q199").Value = 15
$ws.Range("R199").Value = "Hortaliza"
